$wb = $excel.ActiveWorkbook

# Rename the "Include #0" sheet to "Include from LOINC"
$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Name = "Include from LOINC"

# Update the Metadata sheet values (revert of the "Merging 0.1.8 w VitalSigns" merge)
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "0.1.6"
$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2023-05-05T10:50:04-05:00"
$ws.Range("B10").Value = "No display for ContactDetail"
$ws.Range("B11").Value = "No display for ContactDetail"

# Remove the now-obsolete "Jurisdiction" row; everything below shifts up
$ws.Rows.Item(12).Delete()
